$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Rubric": update point values in column D, drop the now-unused
# column E (was blank) and the s="4" style that the cleared-out column E
# cells used to carry.
# ---------------------------------------------------------------------------
$rubric = $wb.Worksheets.Item("Rubric")

$rubric.Range("D6").Value = 8
$rubric.Range("D9").Value = 10
$rubric.Range("D11").Value = 4
$rubric.Range("D12").Value = 4
$rubric.Range("D15").Value = 4

# Strip the old style (numFmtId 0 / fontId 0 "applyFont" xf) off the score
# cells -- they go back to the workbook default style.
$rubric.Range("D6").ClearFormats()
$rubric.Range("D8:D16").ClearFormats()

# Row 7 never got a score in either column -- drop D7 entirely.
$rubric.Range("D7").Clear()

# Column E was only ever empty placeholder cells on this sheet -- remove
# them (and their style) outright.
$rubric.Range("E6:E16").Clear()

# ---------------------------------------------------------------------------
# Sheet "Grade": column B/C only ever held row labels for the "Part 2"
# criteria (rows 10-14); those go away and the real data columns D/E slide
# left into B/C.
# ---------------------------------------------------------------------------
$grade = $wb.Worksheets.Item("Grade")

$grade.Range("B1:C18").Delete(-4159)   # xlShiftToLeft

$grade.Range("B6").Value = 8
$grade.Range("B9").Value = 10
$grade.Range("B11").Value = 4
$grade.Range("B12").Value = 4
$grade.Range("B15").Value = 4

$grade.Range("C6").Value = 8
$grade.Range("C9").Value = 10
$grade.Range("C11").Value = 4
$grade.Range("C12").Value = 4
$grade.Range("C15").Value = 4

# Drop the inherited score-column style, then remove the now fully blank
# row (row 7 no longer has a "possible/score" pair on this sheet).
$grade.Range("B6:C16").ClearFormats()
$grade.Range("B7:C7").Clear()

# ---------------------------------------------------------------------------
# Leave the selection where the editor actually left it on each sheet.
# ---------------------------------------------------------------------------
[void]$rubric.Range("D6:D16").Select()

$grade.Activate()
[void]$grade.Range("A3").Select()
